# chore: update Sheets via scheduled runner
# Refresh the profit-tracking columns (currentAveragePrice / *NQ / *HQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the rows whose
# underlying market data changed, across the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# Leve-profit sheets. A couple of rows flip a profit/loss cell between
# "blank" (no data) and a concrete number, so those use $null to clear the
# cell entirely rather than writing a literal 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 307.5
$ws.Range("I2").Value = 222
$ws.Range("K2").Value = 222
$ws.Range("M2").Value = -109
$ws.Range("H6").Value = 1019.875
$ws.Range("I6").Value = 1102.1428
$ws.Range("K6").Value = 3306.4284
$ws.Range("M6").Value = -3194.4284
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("H17").Value = 1695.8
$ws.Range("J17").Value = 1305.5
$ws.Range("L17").Value = 3916.5
$ws.Range("N17").Value = -4252.5
$ws.Range("H19").Value = 2742.5
$ws.Range("J19").Value = 2993.3333
$ws.Range("L19").Value = 2993.3333
$ws.Range("N19").Value = -3343.3333
$ws.Range("H74").Value = 4859.6
$ws.Range("I74").Value = 4859.6
$ws.Range("K74").Value = 4859.6
$ws.Range("M74").Value = -3923.6
$ws.Range("H77").Value = 4859.6
$ws.Range("I77").Value = 4859.6
$ws.Range("K77").Value = 24298
$ws.Range("M77").Value = -19618
$ws.Range("H121").Value = 1160.7693
$ws.Range("J121").Value = 1174.1666
$ws.Range("L121").Value = 3522.4998
$ws.Range("N121").Value = -7016.4998
$ws.Range("H138").Value = 3176.1667
$ws.Range("I138").Value = 3519
$ws.Range("J138").Value = 2974.5
$ws.Range("K138").Value = 10557
$ws.Range("L138").Value = 8923.5
$ws.Range("M138").Value = -5417
$ws.Range("N138").Value = -19203.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 896753.0600000001
$ws.Range("I2").Value = 1225893.9
$ws.Range("J2").Value = 3371
$ws.Range("K2").Value = 1225893.9
$ws.Range("L2").Value = 3371
$ws.Range("M2").Value = -1225780.9
$ws.Range("N2").Value = -3597
$ws.Range("H32").Value = 969.21
$ws.Range("I32").Value = 896.5789
$ws.Range("K32").Value = 896.5789
$ws.Range("M32").Value = -609.5789
$ws.Range("H102").Value = 1441.2858
$ws.Range("I102").Value = 1441.2858
$ws.Range("K102").Value = 1441.2858
$ws.Range("M102").Value = 180.7141999999999
$ws.Range("H116").Value = 896753.0600000001
$ws.Range("I116").Value = 1225893.9
$ws.Range("J116").Value = 3371
$ws.Range("K116").Value = 1225893.9
$ws.Range("L116").Value = 3371
$ws.Range("M116").Value = -1223599.9
$ws.Range("N116").Value = -7959
$ws.Range("H122").Value = 1411.5454
$ws.Range("I122").Value = 1492
$ws.Range("J122").Value = 607
$ws.Range("K122").Value = 4476
$ws.Range("L122").Value = 1821
$ws.Range("M122").Value = -2026
$ws.Range("N122").Value = -6721

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 896753.0600000001
$ws.Range("I3").Value = 1225893.9
$ws.Range("J3").Value = 3371
$ws.Range("K3").Value = 1225893.9
$ws.Range("L3").Value = 3371
$ws.Range("M3").Value = -1225779.9
$ws.Range("N3").Value = -3599
$ws.Range("H20").Value = 1912.8889
$ws.Range("I20").Value = 1634.8
$ws.Range("K20").Value = 1634.8
$ws.Range("M20").Value = -1387.8
$ws.Range("H22").Value = 112.5
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -371
$ws.Range("H75").Value = 3000
$ws.Range("I75").Value = 3000
$ws.Range("K75").Value = 3000
$ws.Range("M75").Value = -2064
$ws.Range("H78").Value = 3000
$ws.Range("I78").Value = 3000
$ws.Range("K78").Value = 9000
$ws.Range("M78").Value = -4320
$ws.Range("H86").Value = 2471.0625
$ws.Range("J86").Value = 2515.5557
$ws.Range("L86").Value = 2515.5557
$ws.Range("N86").Value = -4761.5557
$ws.Range("H89").Value = 2471.0625
$ws.Range("J89").Value = 2515.5557
$ws.Range("L89").Value = 12577.7785
$ws.Range("N89").Value = -23809.7785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 629.75
$ws.Range("I16").Value = 629.75
$ws.Range("K16").Value = 629.75
$ws.Range("M16").Value = -342.75
$ws.Range("H31").Value = 2082.2046
$ws.Range("I31").Value = 1176.7222
$ws.Range("J31").Value = 2709.077
$ws.Range("K31").Value = 1176.7222
$ws.Range("L31").Value = 2709.077
$ws.Range("M31").Value = -881.7221999999999
$ws.Range("N31").Value = -3299.077
$ws.Range("H34").Value = 2082.2046
$ws.Range("I34").Value = 1176.7222
$ws.Range("J34").Value = 2709.077
$ws.Range("K34").Value = 1176.7222
$ws.Range("L34").Value = 2709.077
$ws.Range("M34").Value = -974.7221999999999
$ws.Range("N34").Value = -3113.077
$ws.Range("H105").Value = 964.6667
$ws.Range("I105").Value = 964.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 964.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 782.3333
$ws.Range("N105").Value = $null
$ws.Range("H113").Value = 629.75
$ws.Range("I113").Value = 629.75
$ws.Range("K113").Value = 629.75
$ws.Range("M113").Value = 1540.25
$ws.Range("H122").Value = 1576.8235
$ws.Range("J122").Value = 2089
$ws.Range("L122").Value = 6267
$ws.Range("N122").Value = -11167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 26950.5
$ws.Range("I55").Value = 100004
$ws.Range("J55").Value = 2599.3333
$ws.Range("K55").Value = 300012
$ws.Range("L55").Value = 7797.999899999999
$ws.Range("M55").Value = -299835
$ws.Range("N55").Value = -8151.999899999999
$ws.Range("H75").Value = 26636
$ws.Range("I75").Value = 600
$ws.Range("J75").Value = 43993.332
$ws.Range("K75").Value = 1800
$ws.Range("L75").Value = 131979.996
$ws.Range("M75").Value = -802
$ws.Range("N75").Value = -133975.996
$ws.Range("H78").Value = 26636
$ws.Range("I78").Value = 600
$ws.Range("J78").Value = 43993.332
$ws.Range("K78").Value = 5400
$ws.Range("L78").Value = 395939.988
$ws.Range("M78").Value = -408
$ws.Range("N78").Value = -405923.988
$ws.Range("H107").Value = 2191.279
$ws.Range("J107").Value = 2246.0833
$ws.Range("L107").Value = 6738.249899999999
$ws.Range("N107").Value = -10578.2499
$ws.Range("H121").Value = 507.5
$ws.Range("I121").Value = 361.25
$ws.Range("J121").Value = 800
$ws.Range("K121").Value = 1083.75
$ws.Range("L121").Value = 2400
$ws.Range("M121").Value = 226.25
$ws.Range("N121").Value = -5020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3160
$ws.Range("I100").Value = 1340
$ws.Range("J100").Value = 4980
$ws.Range("K100").Value = 1340
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -799
$ws.Range("N100").Value = -6062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 525.3333
$ws.Range("I100").Value = 640.5
$ws.Range("J100").Value = 295
$ws.Range("K100").Value = 1281
$ws.Range("L100").Value = 590
$ws.Range("M100").Value = -740
$ws.Range("N100").Value = -1672
$ws.Range("H132").Value = 2471.2632
$ws.Range("I132").Value = 973.25
$ws.Range("K132").Value = 2919.75
$ws.Range("M132").Value = -389.75
